# issue #5: stock data from json to db
# Add "category", "source_file" and "index" columns to the 股票 (stock) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票 (stock) sheet

# Insert a new "category" column before the existing "date" column (I),
# shifting date/legislator_name/legislator_id one column to the right
# (I->J, J->K, K->L). Inserting (rather than just writing new cells)
# carries the header/body styles (s="1" / s="2") along automatically.
$ws.Range("I1:I6").Insert(-4161)

# Insert two more new columns for "source_file" and "index" right after
# the (now shifted) legislator_id column (L), i.e. at M and N.
$ws.Range("M1:N6").Insert(-4161)

# --- Header row ---
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data rows ---
$ws.Range("I2:I6").Value = "normal"
$ws.Range("M2:M6").Value = "tmp81dc1"

$ws.Range("N2").Value = 73
$ws.Range("N3").Value = 74
$ws.Range("N4").Value = 75
$ws.Range("N5").Value = 76
$ws.Range("N6").Value = 77
